# EnumData.xlsx edit:
#  - drop the MapData sheet entirely
#  - rename TestData -> EnumData
#  - add a new "Yellow" ColorType enum row (row 6) in the same style as the
#    existing ColorType rows (3-5)
#  - leave the active selection on D9, matching the saved view state

$wb = $excel.ActiveWorkbook

# Delete the MapData sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("MapData").Delete()
$excel.DisplayAlerts = $true

# Rename TestData -> EnumData
$ws = $wb.Worksheets.Item("TestData")
$ws.Name = "EnumData"

# Row 6 was a blank spacer row; turn it into another ColorType entry.
# Copy A3's formatting (ColorType label style) onto A6 before writing values
# so the new row visually matches rows 3-5.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "ColorType"
$ws.Range("B6").Value = "Yellow"

# Match the saved selection/view state
$ws.Range("D9").Select() | Out-Null
